$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A11").Value = "AV."
$ws.Range("A13").Value = "BA."
$ws.Range("A18").Value = "BP."
$ws.Range("A20").Value = "BT.A"
$ws.Range("A49").Value = "ICG"
$ws.Range("A52").Value = "JD."
$ws.Range("A62").Value = "NG."
$ws.Range("A76").Value = "RR."
$ws.Range("A87").Value = "SN."
$ws.Range("A93").Value = "TW."
$ws.Range("A96").Value = "UU."
